# ---------------------------------------------------------------------------
# Emulacion cambios en control
#
# The authored change (per the commit's xml diff) touches only two spots in
# word/document.xml's single paragraph:
#   1. The drawing run's rPr drops the east-Asian language tag
#      (<w:lang w:eastAsia="es-CO"/>).
#   2. The following run, which holds four literal space characters, is
#      removed entirely.
# Everything else in the paragraph (the drawing itself, the bookmark, the
# sectPr) stays untouched.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument
$targetPara = $d.Paragraphs(1).Range

# Pull the paragraph's own OOXML so we edit the exact content that is
# actually there instead of re-typing the (large) drawing block by hand.
$full = $targetPara.WordOpenXML
$null = $full -match '(?s)<w:p[ >].*?</w:p>'
$paraXml = $Matches[0]

# 1) Remove the <w:lang w:eastAsia="..."/> element from the run properties.
$paraXml = $paraXml -replace '<w:lang w:eastAsia="[^"]*"/>', ''

# 2) Remove the run that contains only whitespace text (the four spaces).
$paraXml = $paraXml -replace '<w:r[^>]*><w:t[^>]*>\s*</w:t></w:r>', ''

# InsertXML needs the fragment to carry its own namespace declarations
# since it is no longer nested inside the original <w:document> element.
$nsDecls = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
           'xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" ' +
           'xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" ' +
           'xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" ' +
           'xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" ' +
           'xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" ' +
           'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'
$paraXml = $paraXml -replace '^<w:p', "<w:p $nsDecls"

$targetPara.InsertXML($paraXml)
